# Updates cryptos list data (prices and volume percentages) per commit:
# "Updated cryptos list on Wed Nov 20 05:56:50 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force pure-numeric-looking Price values to remain stored as text,
# matching the source data which represents all Price/Volume cells as strings.
$forceTextCells = @(
    "D5",
    "D6",
    "D8",
    "D11",
    "D13",
    "D15",
    "D16",
    "D19",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D42",
    "D43",
    "D44",
    "D46",
    "D47",
    "D48",
    "D49"
)
foreach ($cellRef in $forceTextCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "92.483.29"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "3.109.09"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("D5").Value = "234.59"
$ws.Range("E5").Value = "  -2.73%  "
$ws.Range("D6").Value = "613.12"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  -2.62%  "
$ws.Range("D8").Value = "0.390"
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "3.106.49"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "0.782"
$ws.Range("E11").Value = "  +3.75%  "
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("D13").Value = "0.0000244"
$ws.Range("E13").Value = "  -4.38%  "
$ws.Range("D14").Value = "92.195.00"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "33.87"
$ws.Range("E15").Value = "  -4.11%  "
$ws.Range("D16").Value = "5.42"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("D17").Value = "3.686.71"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").Value = "3.065.80"
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("D19").Value = "3.80"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").Value = "5.82"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").Value = "0.0000205"
$ws.Range("E22").Value = "  +1.17%  "
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "438.94"
$ws.Range("E24").Value = "  -4.48%  "
$ws.Range("D25").Value = "5.58"
$ws.Range("E25").Value = "  -5.99%  "
$ws.Range("D26").Value = "85.25"
$ws.Range("E26").Value = "  -4.40%  "
$ws.Range("D27").Value = "11.51"
$ws.Range("E27").Value = "  -2.35%  "
$ws.Range("D28").Value = "3.267.44"
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "0.177"
$ws.Range("E30").Value = "  +5.41%  "
$ws.Range("D31").Value = "0.230"
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -20.66%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "1.04"
$ws.Range("E33").Value = "  -31.60%  "
$ws.Range("D34").Value = "9.18"
$ws.Range("E34").Value = "  -2.24%  "
$ws.Range("D35").Value = "8.02"
$ws.Range("E35").Value = "  +7.37%  "
$ws.Range("D36").Value = "0.157"
$ws.Range("E36").Value = "  -9.86%  "
$ws.Range("D37").Value = "25.86"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").Value = "3.98"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").Value = "1.90"
$ws.Range("E39").Value = "  -3.62%  "
$ws.Range("D40").Value = "23.87"
$ws.Range("E40").Value = "  +7.67%  "
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("D42").Value = "467.00"
$ws.Range("E42").Value = "  -5.29%  "
$ws.Range("D43").Value = "0.431"
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("D44").Value = "3.27"
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "159.87"
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("D47").Value = "0.684"
$ws.Range("E47").Value = "  -3.55%  "
$ws.Range("D48").Value = "1.84"
$ws.Range("E48").Value = "  -4.86%  "
$ws.Range("D49").Value = "0.0330"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("E51").Value = "  -0.58%  "
